# Note_sheet_CH.xlsx - "Fixed date-time info" edit.
#
# Column B ("Time") held the clock times as plain text strings like
# "5:15PM" (stored as shared-string cells). That meant Excel could not
# treat them as real times (no sorting/arithmetic/formatting as a time).
# This script rewrites every Time cell on Sheet1 as a genuine numeric
# Excel time serial (fraction of a 24h day) and applies a proper custom
# time number format to the column, which is the actual "date-time info"
# fix referenced by the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Original text that used to live in each B-column cell (row -> "h:mmAM/PM").
$timeByRow = @{
    2  = "5:15PM"
    3  = "5:15PM"
    4  = "5:00PM"
    5  = "5:00PM"
    6  = "5:00PM"
    7  = "6:00PM"
    8  = "3:15PM"
    9  = "3:15PM"
    10 = "4:15PM"
    11 = "4:15PM"
    12 = "3:30PM"
    13 = "3:30PM"
    14 = "4:30PM"
    15 = "3:30PM"
    16 = "5:30PM"
    17 = "6:30PM"
    18 = "6:30PM"
    19 = "6:30PM"
    20 = "12:00PM"
    21 = "12:00PM"
    22 = "12:00PM"
    23 = "12:00PM"
    24 = "2:00PM"
    25 = "5:15PM"
    26 = "5:15PM"
    27 = "6:15PM"
    28 = "5:00PM"
    29 = "5:00PM"
    30 = "6:00PM"
    31 = "6:00PM"
    32 = "3:15PM"
    33 = "4:15PM"
    34 = "4:15PM"
    35 = "5:15PM"
    36 = "3:30PM"
    37 = "4:30PM"
    38 = "4:30PM"
    39 = "5:00PM"
    40 = "5:00PM"
    41 = "6:00PM"
    42 = "6:00PM"
    43 = "5:00PM"
    44 = "5:00PM"
    45 = "6:00PM"
    46 = "3:15PM"
    47 = "4:15PM"
    48 = "4:15pm"
    49 = "5:15PM"
    50 = "5:15PM"
    51 = "3:30PM"
    52 = "3:30pm"
}

# Parse a "h:mmAM/PM" (case-insensitive) clock string into an Excel time
# serial - the fraction of a day it represents - without relying on
# locale-sensitive DateTime parsing.
function ConvertTo-DayFraction($clockText) {
    $suffix = $clockText.Substring($clockText.Length - 2, 2).ToUpper()
    $hourMinute = $clockText.Substring(0, $clockText.Length - 2)
    $segments = $hourMinute.Split(":")
    $hour = [int]$segments[0]
    $minute = [int]$segments[1]
    if ($suffix -eq "PM" -and $hour -ne 12) { $hour = $hour + 12 }
    if ($suffix -eq "AM" -and $hour -eq 12) { $hour = 0 }
    return ($hour * 3600 + $minute * 60) / 86400
}

foreach ($row in $timeByRow.Keys) {
    $cell = $ws.Range("B" + $row)
    $cell.Value2 = ConvertTo-DayFraction $timeByRow[$row]
}

# Give the now-numeric Time column a real time display format (was the
# generic h:mm AM/PM built-in before; this matches what Excel wrote when
# the format was picked from Format Cells > Time > the long-time locale
# format).
$ws.Range("B2:B52").NumberFormat = "[$-F400]h:mm:ss\ AM/PM"

# Column got visibly wider once it held real formatted times - autofit it.
$ws.Columns("B").AutoFit()

# Printer/page setup was touched (explicit Portrait orientation) and the
# last active selection moved to C50.
$ws.PageSetup.Orientation = [Microsoft.Office.Interop.Excel.XlPageOrientation]::xlPortrait
$ws.Range("C50").Select() | Out-Null
